$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lương")

$ws.Range("B12").Value = 20
$ws.Range("B13").Value = 700000
$ws.Range("B14").Value = 2142857.142857143
$ws.Range("B33").Value = 1512857.142857143
$ws.Range("B35").Value = 1512857.142857143
